$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (approximated via ColumnWidth COM property; the runtime
# quantizes the stored XML width to the nearest 1/6 after adding 5/6, so
# these inputs are chosen to land as close as possible to the target widths)
$ws.Columns.Item(1).ColumnWidth = 12.366666666666665
$ws.Columns.Item(2).ColumnWidth = 17.166666666666668
$ws.Columns.Item(3).ColumnWidth = 7.566666666666667
$ws.Columns.Item(4).ColumnWidth = 8.766666666666666
$ws.Columns.Item(5).ColumnWidth = 8.766666666666666
$ws.Columns.Item(6).ColumnWidth = 7.566666666666667
$ws.Columns.Item(7).ColumnWidth = 5.166666666666667
$ws.Columns.Item(8).ColumnWidth = 14.766666666666666
$ws.Columns.Item(9).ColumnWidth = 6.366666666666667
$ws.Columns.Item(10).ColumnWidth = 59.166666666666664

# Convert fractional values (0-1) to percentage values (0-100), rounded to 2dp
$ws.Range("B2").Value = 11.69
$ws.Range("C2").Value = 34.42
$ws.Range("D2").Value = 27.27
$ws.Range("E2").Value = 9.74
$ws.Range("F2").Value = 2.6
$ws.Range("G2").Value = 7.79
$ws.Range("H2").Value = 6.49
$ws.Range("I2").Value = 0.0
$ws.Range("J2").Value = 100.0
$ws.Range("B3").Value = 0.0
$ws.Range("C3").Value = 46.99
$ws.Range("D3").Value = 36.75
$ws.Range("E3").Value = 5.42
$ws.Range("F3").Value = 7.23
$ws.Range("G3").Value = 3.61
$ws.Range("H3").Value = 0.0
$ws.Range("I3").Value = 0.0
$ws.Range("J3").Value = 100.0
$ws.Range("B4").Value = 17.31
$ws.Range("C4").Value = 23.08
$ws.Range("D4").Value = 25.64
$ws.Range("E4").Value = 12.18
$ws.Range("F4").Value = 5.77
$ws.Range("G4").Value = 8.33
$ws.Range("H4").Value = 7.69
$ws.Range("I4").Value = 0.0
$ws.Range("J4").Value = 100.0
$ws.Range("B5").Value = 0.0
$ws.Range("C5").Value = 49.68
$ws.Range("D5").Value = 17.2
$ws.Range("E5").Value = 22.29
$ws.Range("F5").Value = 6.37
$ws.Range("G5").Value = 4.46
$ws.Range("H5").Value = 0.0
$ws.Range("I5").Value = 0.0
$ws.Range("J5").Value = 100.0
$ws.Range("B6").Value = 0.0
$ws.Range("C6").Value = 29.52
$ws.Range("D6").Value = 24.7
$ws.Range("E6").Value = 19.28
$ws.Range("F6").Value = 9.04
$ws.Range("G6").Value = 13.86
$ws.Range("H6").Value = 3.61
$ws.Range("I6").Value = 0.0
$ws.Range("J6").Value = 100.0
$ws.Range("B7").Value = 0.0
$ws.Range("C7").Value = 22.53
$ws.Range("D7").Value = 36.26
$ws.Range("E7").Value = 16.48
$ws.Range("F7").Value = 9.89
$ws.Range("G7").Value = 12.64
$ws.Range("H7").Value = 2.2
$ws.Range("I7").Value = 0.0
$ws.Range("J7").Value = 100.0
$ws.Range("B8").Value = 15.13
$ws.Range("C8").Value = 36.84
$ws.Range("D8").Value = 19.74
$ws.Range("E8").Value = 9.87
$ws.Range("F8").Value = 3.95
$ws.Range("G8").Value = 9.21
$ws.Range("H8").Value = 5.26
$ws.Range("I8").Value = 0.0
$ws.Range("J8").Value = 100.0
$ws.Range("B9").Value = 20.99
$ws.Range("C9").Value = 8.02
$ws.Range("D9").Value = 20.37
$ws.Range("E9").Value = 22.22
$ws.Range("F9").Value = 8.64
$ws.Range("G9").Value = 12.35
$ws.Range("H9").Value = 7.41
$ws.Range("I9").Value = 0.0
$ws.Range("J9").Value = 100.0
$ws.Range("B10").Value = 0.0
$ws.Range("C10").Value = 14.2
$ws.Range("D10").Value = 51.85
$ws.Range("E10").Value = 22.84
$ws.Range("F10").Value = 8.64
$ws.Range("G10").Value = 2.47
$ws.Range("H10").Value = 0.0
$ws.Range("I10").Value = 0.0
$ws.Range("J10").Value = 100.0
$ws.Range("B11").Value = 1.81
$ws.Range("C11").Value = 22.29
$ws.Range("D11").Value = 25.3
$ws.Range("E11").Value = 21.69
$ws.Range("F11").Value = 12.65
$ws.Range("G11").Value = 12.65
$ws.Range("H11").Value = 3.61
$ws.Range("I11").Value = 0.0
$ws.Range("J11").Value = 100.0
